# Applies the "Experimental setup made!" edit to the meeting-minutes document.
#
# Word normally keeps adjacent runs that share identical formatting split
# apart (e.g. "-" and " How well does our robot clean?" in two <w:r>'s).
# When the user retypes / edits such a line in real Word, the runs that make
# up the edited line get coalesced into a single run. We reproduce that by
# replacing the *whole* line's Range.Text in one shot using a helper.

$d = $word.ActiveDocument
$vt = [char]11          # w:br (line break) character inside paragraph text
$MARKER = "@@EDITMARK@@"

function Merge-BreakRun($doc, $paraIndex, $breakOccurrence, $newTail) {
    # Rewrites the text from the N-th line-break (1-based) inside paragraph
    # $paraIndex through the following line-break (or paragraph end) so
    # that it becomes a single run "<br/>$newTail". Word/this host only
    # restructures runs when the assigned text genuinely differs from what
    # is already there, so we append a throwaway marker to force a real
    # change, then delete the marker in a second, separate pass.
    $p = $doc.Paragraphs.Item($paraIndex)
    $t = $p.Range.Text

    $positions = New-Object System.Collections.ArrayList
    for ($j = 0; $j -lt $t.Length; $j++) {
        if ([int][char]$t[$j] -eq 11) {
            [void]$positions.Add($j)
        }
    }

    $idx = [int]$positions[$breakOccurrence - 1]
    $pStart = $p.Range.Start
    $absBreak = $pStart + $idx

    if ($breakOccurrence -lt $positions.Count) {
        $endIdx = [int]$positions[$breakOccurrence]
    } else {
        $endIdx = $t.Length
    }
    $absEnd = $pStart + $endIdx

    $r = $doc.Range($absBreak, $absEnd)
    $r.Text = $vt + $newTail + $MARKER

    $p2 = $doc.Paragraphs.Item($paraIndex)
    $t2 = $p2.Range.Text
    $mpos = $t2.IndexOf($MARKER)
    $absM = $p2.Range.Start + $mpos
    $rm = $doc.Range($absM, $absM + $MARKER.Length)
    $rm.Delete()
}

function Append-BreakLine($doc, $paraIndex, $newLine) {
    # Appends "<br/>$newLine" at the end of paragraph $paraIndex (before the
    # paragraph mark).
    $p = $doc.Paragraphs.Item($paraIndex)
    $insPos = $p.Range.End - 1
    $ins = $doc.Range($insPos, $insPos)
    $ins.InsertAfter($vt + $newLine)
}

# --- "04-05-2022, Wednesday week 2" section -------------------------------

# "Questions to answer:" paragraph -> merge "-" + " How well..." runs
Merge-BreakRun $d 9 2 "- How well does our robot clean?"

# "Assumptions:" paragraph -> merge "- There's a" + " load station..." runs
Merge-BreakRun $d 10 2 "- There's a load station at starting point"

# "Three different space settings:" paragraph:
#  - merge the break before "Obstacles versus no obstacles [2]" into its run
#  - append the new "- Death tiles yes/no [2]" line
Merge-BreakRun $d 11 3 "- Obstacles versus no obstacles [2]"
Append-BreakLine $d 11 "- Death tiles yes/no [2]"

# "===> Resulting in 2³ = 8 different spaces to clean" paragraph becomes
# "===> Resulting in 2^4 = 16 different spaces to clean", with the
# _GoBack bookmark now sitting right after "16".
$p12 = $d.Paragraphs.Item(12)
$p12Start = $p12.Range.Start
$p12End = $p12.Range.End
$rFull = $d.Range($p12Start, $p12End - 1)
$rFull.Text = "===> Resulting in 2^4 = 16 different spaces to clean"

$t12 = $d.Paragraphs.Item(12).Range.Text
$markerPos = $t12.IndexOf("16") + 2
$bmPos = $p12Start + $markerPos
$bmRange = $d.Range($bmPos, $bmPos)
# Re-adding a bookmark with an existing name moves it, automatically
# removing it from its old location near "- Number of robots?".
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# "Evaluation metrics:" paragraph -> merge split runs on two of its lines
Merge-BreakRun $d 13 2 "- % efficiency (#tiles cleaned/#tiles visited)"
Merge-BreakRun $d 13 4 "- Number of iterations?"

# "Other parameters to vary:" paragraph -> merge final "- Number of robots?"
# line (the bookmark that used to sit here was relocated above).
Merge-BreakRun $d 14 2 "- Number of robots?"
